# e2e flow with file downloading: populate the result sheet with the
# scraped product URLs (one per row in column A), each as a hyperlink
# back to the product page with the option id as the in-page anchor.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 already has a hyperlink from the previous run; replace it ---
# --- in place (drop the old link, add the new one) so it keeps rId1. ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "https://makeup.com.ua/product/8383/", "/option/463699/")
$ws.Range("A1").Value = "https://makeup.com.ua/product/8383/#/option/463699/"

# --- Rows 2-14: new rows, each with its own hyperlink. ---
$ws.Hyperlinks.Add($ws.Range("A2"), "https://makeup.com.ua/product/8383/", "/option/463697/", "", "https://makeup.com.ua/product/8383/ - /option/463697/")
$ws.Range("A2").Value = "https://makeup.com.ua/product/8383/#/option/463697/"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://makeup.com.ua/product/3173/", "/option/401569/", "", "https://makeup.com.ua/product/3173/ - /option/401569/")
$ws.Range("A3").Value = "https://makeup.com.ua/product/3173/#/option/401569/"
$ws.Hyperlinks.Add($ws.Range("A4"), "https://makeup.com.ua/product/3173/", "/option/401571/", "", "https://makeup.com.ua/product/3173/ - /option/401571/")
$ws.Range("A4").Value = "https://makeup.com.ua/product/3173/#/option/401571/"
$ws.Hyperlinks.Add($ws.Range("A5"), "https://makeup.com.ua/product/3173/", "/option/401567/", "", "https://makeup.com.ua/product/3173/ - /option/401567/")
$ws.Range("A5").Value = "https://makeup.com.ua/product/3173/#/option/401567/"
$ws.Hyperlinks.Add($ws.Range("A6"), "https://makeup.com.ua/product/583163/", "/option/1474529/", "", "https://makeup.com.ua/product/583163/ - /option/1474529/")
$ws.Range("A6").Value = "https://makeup.com.ua/product/583163/#/option/1474529/"
$ws.Hyperlinks.Add($ws.Range("A7"), "https://makeup.com.ua/product/215167/", "/option/471781/", "", "https://makeup.com.ua/product/215167/ - /option/471781/")
$ws.Range("A7").Value = "https://makeup.com.ua/product/215167/#/option/471781/"
$ws.Hyperlinks.Add($ws.Range("A8"), "https://makeup.com.ua/product/11180/", "/option/472003/", "", "https://makeup.com.ua/product/11180/ - /option/472003/")
$ws.Range("A8").Value = "https://makeup.com.ua/product/11180/#/option/472003/"
$ws.Hyperlinks.Add($ws.Range("A9"), "https://makeup.com.ua/product/11180/", "/option/471999/", "", "https://makeup.com.ua/product/11180/ - /option/471999/")
$ws.Range("A9").Value = "https://makeup.com.ua/product/11180/#/option/471999/"
$ws.Hyperlinks.Add($ws.Range("A10"), "https://rozetka.com.ua/antonio_banderas_8411061636275/p2215142/")
$ws.Hyperlinks.Add($ws.Range("A11"), "https://rozetka.com.ua/versace_img479/p11514465/")
$ws.Hyperlinks.Add($ws.Range("A12"), "https://rozetka.com.ua/versace_8011003996025/p57030999/")
$ws.Hyperlinks.Add($ws.Range("A13"), "https://rozetka.com.ua/creed_3508440505118/p69683708/")
$ws.Hyperlinks.Add($ws.Range("A14"), "https://rozetka.com.ua/creed_3508441001114/p69673676/")

# --- All 14 link cells share the same built-in "Hyperlink" cell style. ---
$ws.Range("A1:A14").Style = "Hyperlink"

# --- Widen column A to fit the longer URLs, and leave the selection on ---
# --- the last populated cell, matching where the author left off. ---
$ws.Columns.Item(1).ColumnWidth = 66.66666667
$ws.Range("A14").Select()
